# Add lifelines to command result
# Slide 1 ("TaskList" sequence) gains a dashed blue lifeline connector
# (id 52 "Straight Connector 51") and Slide 2 ("TaskListFind" sequence)
# gains a matching one (id 63 "Straight Connector 62"), both drawn the
# same way the existing lifeline connectors on the slides are built:
# a flipped <a:prstGeom prst="line"> connector, 1.5pt weight, solid blue
# (0070C0), system-dash style.

$p = $ppt.ActivePresentation

function Add-Lifeline {
    param($Slide, $ShapeName, $OffX, $OffY, $ExtCx, $ExtCy)

    $emuPerPt = 12700.0
    $beginX = ($OffX + $ExtCx) / $emuPerPt
    $beginY = $OffY / $emuPerPt
    $endX = $OffX / $emuPerPt
    $endY = ($OffY + $ExtCy) / $emuPerPt

    $cxn = $Slide.Shapes.AddLine($beginX, $beginY, $endX, $endY)
    $cxn.Name = $ShapeName
    $cxn.HorizontalFlip = -1

    $cxn.Line.Visible = -1
    $cxn.Line.Weight = 1.5
    $cxn.Line.ForeColor.RGB = 0xC07000
    $cxn.Line.DashStyle = 9

    $cxn.Shadow.Visible = $false

    return $cxn
}

$slide1 = $p.Slides.Item(1)
Add-Lifeline $slide1 "Straight Connector 51" 8077200 5534171 4512 1030836 | Out-Null

$slide2 = $p.Slides.Item(2)
Add-Lifeline $slide2 "Straight Connector 62" 8846012 5522364 4512 1030836 | Out-Null
